# Landscaping Data.xlsx - "Add files via upload"
# Appends 7 new observation rows (261-267) for 2025-06-16 (serial 45824),
# extending the data table that previously ended at row 260.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Stamp out 7 fresh rows with the same formatting as the last
#        existing data row (this preserves the date number format on
#        column A and the default/general formatting on the rest). ---
for ($r = 261; $r -le 267; $r++) {
    $ws.Range("A260:T260").Copy($ws.Range("A" + $r + ":T" + $r))
}

# --- 2. New row data (Date, Plant_Type, Plant_Size, Low, High, Rain,
#        Growth, Pruned, Quadrant, Shade, UV, Humidity, Dew_Point,
#        Pressure, Wind_Gust, Cloud_Cover, Visibility, AQI, Pollen). ---
$rows = @(
    @{ Row=261; B="Flowering";     C="Large";  D=67; E=74; G=0.49; H=0.1;  I="No"; J=2; K="Neutral"; L=4; M=0.85; N=69; O=30.01; P=5; Q=0.81; R=9.9; S=22; T=39 },
    @{ Row=262; B="Nonflowering";  C="Medium"; D=67; E=74; G=0.49; H=0.1;  I="No"; J=3; K="Neutral"; L=4; M=0.85; N=69; O=30.01; P=5; Q=0.81; R=9.9; S=22; T=39 },
    @{ Row=263; B="Nonflowering";  C="Small";  D=67; E=74; G=0.49; H=0.2;  I="No"; J=3; K="Neutral"; L=4; M=0.85; N=69; O=30.01; P=5; Q=0.81; R=9.9; S=22; T=39 },
    @{ Row=264; B="Nonflowering";  C="Medium"; D=67; E=74; G=0.49; H=0.3;  I="No"; J=3; K="Bright";  L=4; M=0.85; N=69; O=30.01; P=5; Q=0.81; R=9.9; S=22; T=39 },
    @{ Row=265; B="Nonflowering";  C="Medium"; D=67; E=74; G=0.49; H=0.25; I="No"; J=3; K="Bright";  L=4; M=0.85; N=69; O=30.01; P=5; Q=0.81; R=9.9; S=22; T=39 },
    @{ Row=266; B="Nonflowering";  C="Large";  D=67; E=74; G=0.49; H=0.25; I="No"; J=4; K="Neutral"; L=4; M=0.85; N=69; O=30.01; P=5; Q=0.81; R=9.9; S=22; T=39 },
    @{ Row=267; B="Tree";          C="Medium"; D=67; E=74; G=0.49; H=1.1;  I="No"; J=1; K="Dark";    L=4; M=0.85; N=69; O=30.01; P=5; Q=0.81; R=9.9; S=22; T=39 }
)

foreach ($rowData in $rows) {
    $r = $rowData.Row

    $ws.Cells.Item($r, 1).Value2 = 45824
    $ws.Cells.Item($r, 2).Value = $rowData.B
    $ws.Cells.Item($r, 3).Value = $rowData.C
    $ws.Cells.Item($r, 4).Value2 = $rowData.D
    $ws.Cells.Item($r, 5).Value2 = $rowData.E
    $ws.Cells.Item($r, 6).Formula = "=ABS(D" + $r + "-E" + $r + ")"
    $ws.Cells.Item($r, 7).Value2 = $rowData.G
    $ws.Cells.Item($r, 8).Value2 = $rowData.H
    $ws.Cells.Item($r, 9).Value = $rowData.I
    $ws.Cells.Item($r, 10).Value2 = $rowData.J
    $ws.Cells.Item($r, 11).Value = $rowData.K
    $ws.Cells.Item($r, 12).Value2 = $rowData.L
    $ws.Cells.Item($r, 13).Value2 = $rowData.M
    $ws.Cells.Item($r, 14).Value2 = $rowData.N
    $ws.Cells.Item($r, 15).Value2 = $rowData.O
    $ws.Cells.Item($r, 16).Value2 = $rowData.P
    $ws.Cells.Item($r, 17).Value2 = $rowData.Q
    $ws.Cells.Item($r, 18).Value2 = $rowData.R
    $ws.Cells.Item($r, 19).Value2 = $rowData.S
    $ws.Cells.Item($r, 20).Value2 = $rowData.T
}
